# Apply the update described by the diff:
#  - Column A (Timestamp) for rows 2..97: every value shifts forward by
#    exactly one day (the fetched day moves forward by one).
#  - Column B (Notified Production (MW)) for rows 2..93: replaced with the
#    new notified production values from the updated ENTSO-E fetch.
#  - Rows 94..97 in column B remain 0 (unchanged by the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp (serial date) values for A2:A97.
$aValues = @(
    45734.01041666666,45734.02083333334,45734.03125,45734.04166666666,45734.05208333334,45734.0625,45734.07291666666,45734.08333333334,45734.09375,45734.10416666666,
    45734.11458333334,45734.125,45734.13541666666,45734.14583333334,45734.15625,45734.16666666666,45734.17708333334,45734.1875,45734.19791666666,45734.20833333334,
    45734.21875,45734.22916666666,45734.23958333334,45734.25,45734.26041666666,45734.27083333334,45734.28125,45734.29166666666,45734.30208333334,45734.3125,
    45734.32291666666,45734.33333333334,45734.34375,45734.35416666666,45734.36458333334,45734.375,45734.38541666666,45734.39583333334,45734.40625,45734.41666666666,
    45734.42708333334,45734.4375,45734.44791666666,45734.45833333334,45734.46875,45734.47916666666,45734.48958333334,45734.5,45734.51041666666,45734.52083333334,
    45734.53125,45734.54166666666,45734.55208333334,45734.5625,45734.57291666666,45734.58333333334,45734.59375,45734.60416666666,45734.61458333334,45734.625,
    45734.63541666666,45734.64583333334,45734.65625,45734.66666666666,45734.67708333334,45734.6875,45734.69791666666,45734.70833333334,45734.71875,45734.72916666666,
    45734.73958333334,45734.75,45734.76041666666,45734.77083333334,45734.78125,45734.79166666666,45734.80208333334,45734.8125,45734.82291666666,45734.83333333334,
    45734.84375,45734.85416666666,45734.86458333334,45734.875,45734.88541666666,45734.89583333334,45734.90625,45734.91666666666,45734.92708333334,45734.9375,
    45734.94791666666,45734.95833333334,45734.96875,45734.97916666666,45734.98958333334,45735
)

# New "Notified Production (MW)" values for B2:B93. Rows 94-97 stay 0.
$bValues = @(
    2579,2585,2603,2625,2672,2673,2667,2665,2684,2684,
    2654,2653,2668,2665,2664,2654,2652,2655,2650,2642,
    2606,2608,2604,2608,2600,2596,2599,2597,2583,2583,
    2584,2585,2554,2554,2557,2558,2517,2517,2516,2516,
    2453,2452,2450,2449,2368,2365,2363,2361,2293,2290,
    2288,2285,2250,2247,2245,2243,2209,2207,2205,2202,
    2104,2101,2097,2094,1909,1899,1889,1880,1580,1572,
    1564,1556,1309,1303,1298,1293,1135,1132,1129,1126,
    1048,1047,1045,1044,1010,1009,1008,1008,993,991,
    990,989
)

$row = 2
foreach ($val in $aValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row++
}

$row = 2
foreach ($val in $bValues) {
    $ws.Cells.Item($row, 2).Value = $val
    $row++
}
